$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "new transfer to uts" - shift every value in B2:B101 by a constant offset
# (new unit-transfer-system calibration constant).
$offset = 33.3184521174685

for ($row = 2; $row -le 101; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = $cell.Value()
    $cell.Value = $old + $offset
}
